# Apply updated parameter values to the "potential_preg_untrt" sheet and
# switch the active sheet / selection to reflect what the user was looking
# at when they saved the workbook.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("potential_preg_untrt")

# Updated probabilities (column C) for rows 9-17.
$ws.Range("C9").Value  = 0.05
$ws.Range("C10").Value = 0.02
$ws.Range("C11").Value = 0.02
$ws.Range("C13").Value = 0.005
$ws.Range("C14").Value = 0.004
$ws.Range("C15").Value = 0.004
$ws.Range("C16").Value = 0.004
$ws.Range("C17").Value = 0.004

# Make this the active sheet and select the range the author had selected.
$ws.Activate()
$ws.Range("C2:C21").Select()
